$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "halosalsa2" username/email pair to "halosalsa3"
$ws.Range("B7").Value = "halosalsa3"
$ws.Range("B8").Value = "halosalsa3@gmail.com"
